$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update existing sheet "o_10": add column E header + new row-2 values
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "o_10"

$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)
$ws1.Range("E1").Value = "evaluator_partial_correctness"

$promptO10 = @'
 Given is the adjacency matrix for a weighted directed graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the least cost path from node A to node P?
   A B C D E F G H I J K L M N O P
 A 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 5 0 5 0 0 5 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 5 0 0 1 0 0 0 0 0 0 0
 F 0 0 0 0 0 0 4 0 0 5 0 0 0 0 0 0
 G 0 0 5 0 0 0 0 0 0 0 1 0 0 0 0 0
 H 0 0 0 3 0 0 4 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 2 0 0 0 0 1 0
 L 0 0 0 0 0 0 0 2 0 0 4 0 0 0 0 2
 M 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 2 0 0 5 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 2
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0

Solution: No possible path from A to P
        
 Given these examples, answer the following quesiton.

what is the least cost path from node A to node P?

   A B C D E F G H I J K L M N O P
 A 0 4 0 0 3 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 5 0 0 0 0 1 0 0 0 0 0 0
 G 0 0 5 0 0 4 0 1 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 1 0 0 0 0 5 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0
 K 0 0 0 0 0 0 5 0 0 5 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 2 0 0 2 0 0 0 0 1
 M 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 4
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
    
'@

$llmO10 = @'
The least cost path from node A to node P is A -> B -> C -> I -> J -> N -> O -> P with a total cost of 12.
'@

$ws1.Range("A2").Value = $promptO10
$ws1.Range("B2").Value = "No possible path from A to P"
$ws1.Range("C2").Value = $llmO10
$ws1.Range("D2").Value = "Wrong"
$ws1.Range("E2").Value = "Output: 0/1"

# ---------------------------------------------------------------------------
# 2) Add sheet "o_20" right after "o_10"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "o_20"

$ws1.Range("A1:E1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)

$ws2.Range("A1").Value = "prompt"
$ws2.Range("B1").Value = "solution"
$ws2.Range("C1").Value = "llm_response"
$ws2.Range("D1").Value = "evaluator_response"
$ws2.Range("E1").Value = "evaluator_partial_correctness"

$promptO20 = @'
 Given is the adjacency matrix for a weighted directed graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the least cost path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 3 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 5 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 3 0 0 0 0 0 4 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 3 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 2 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 2 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 3 0 0 0 1 0 4 0 0 0 5 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 4 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 2 0 0 0 0 0 5 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 3 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 1 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0
Solution: No possible path from A to Y
 Given these examples, answer the following quesiton.
what is the least cost path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 2 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 4 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 4 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 4 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 2 0 0 0 3 0 3 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 2 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 4 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 4 0 0 0 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 3 0 3 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0
    
'@

$llmO20 = @'
The least cost path from node A to node Y is: A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T -> U -> V -> W -> X -> Y. 
The total cost of this path is 36.
'@

$ws2.Range("A2").Value = $promptO20
$ws2.Range("B2").Value = "No possible path from A to Y"
$ws2.Range("C2").Value = $llmO20
$ws2.Range("D2").Value = "Wrong"
$ws2.Range("E2").Value = "0/1"

# ---------------------------------------------------------------------------
# 3) Add sheet "o_20_jumbled" right after "o_20"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "o_20_jumbled"

$ws1.Range("A1:E1").Copy()
$ws3.Range("A1").PasteSpecial(-4122)

$ws3.Range("A1").Value = "prompt"
$ws3.Range("B1").Value = "solution"
$ws3.Range("C1").Value = "llm_response"
$ws3.Range("D1").Value = "evaluator_response"
$ws3.Range("E1").Value = "evaluator_partial_correctness"

$promptJumbled = @'
 Given is the adjacency matrix for a weighted directed graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the least cost path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 3 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 5 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 3 0 0 0 0 0 4 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 3 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 2 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 2 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 3 0 0 0 1 0 4 0 0 0 5 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 4 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 2 0 0 0 0 0 5 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 3 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 1 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0
Solution: No possible path from A to Y
 Given these examples, answer the following quesiton.
what is the least cost path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 5 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 5 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 4 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 3 0 0 0 4 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 4 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 5 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 2 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 3 0 4 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 2 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 2 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 5
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 3 0 3 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 2
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
    
'@

$llmJumbled = @'
The least cost path from node A to node Y is as follows:
A -> B -> G -> H -> I -> J -> N -> S -> X -> Y
Total cost = 1 + 1 + 1 + 3 + 5 + 3 + 3 + 2 + 3 = 22
'@

$ws3.Range("A2").Value = $promptJumbled
$ws3.Range("B2").Value = "No possible path from A to Y"
$ws3.Range("C2").Value = $llmJumbled
$ws3.Range("D2").Value = "Wrong"
$ws3.Range("E2").Value = "Output: 0/1"

# ---------------------------------------------------------------------------
# 4) Restore original active sheet/tab selection
# ---------------------------------------------------------------------------
$ws1.Activate()
